# Auto-generated edit script
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final target values for rows 100..166, columns D,J,K,L,M,O,P
# (a new weekly record was inserted as row 100, shifting the remaining
# historical rows down by one, and the former last row reappears as row 166)
$startRow = 100
$colD = @(44582,44413,44544,44357,44320,44306,44407,44316,44329,44460,44526,44467,44411,44313,44313,44334,44389,44517,44221,44280,44330,44483,44448,44463,44239,44476,44169,44250,44515,44379,44509,44229,44426,44488,44341,44452,44473,44298,44305,44301,44278,44466,44392,44412,44322,44236,44162,44434,44532,44490,44427,44491,44293,44533,44494,44571,44390,44279,44481,44277,44525,44327,44354,44503,44462,44312,44511)
$colJ = @(280,120,280,80,240,150,270,240,80,80,240,240,250,200,240,240,40,40,50,50,240,120,120,240,240,120,250,200,80,240,240,150,40,240,240,80,80,60,40,80,120,90,180,40,40,160,250,80,120,80,120,240,40,240,80,80,240,40,300,40,80,260,40,40,120,40,80)
$colK = @(2500,4500,2500,3000,2500,5500,4500,3000,2500,4500,2500,3500,4500,5500,3000,3000,3000,2500,6000,5000,3000,2500,4000,4000,5000,3000,3000,6000,2500,3000,2500,6000,5000,2500,2500,4500,3000,6000,5500,5000,5000,4000,3000,4500,2500,5000,3000,5000,2500,2500,5000,2500,6000,2500,2500,3000,3000,5000,2500,5000,2500,2500,3000,2500,4000,6000,2500)
$colL = @(3000,4500,2500,3000,2500,5500,5000,3000,2500,4500,2500,3500,4500,6000,3000,3000,3000,2500,6000,5000,3000,2500,4000,4000,5000,3000,3000,6000,2500,3000,2500,6000,5000,2500,2500,4500,3000,6000,5500,5000,5000,4000,3000,4500,2500,6000,3000,5000,2500,2500,5000,2500,6000,2500,2500,3000,3000,5000,2500,5000,2500,3000,3000,2500,4000,6000,2500)
$colM = @(2750,4500,2500,3000,2500,5500,4778,3000,2500,4500,2500,3500,4500,5750,3000,3000,3000,2500,6000,5000,3000,2500,4000,4000,5000,3000,3000,6000,2500,3000,2500,6000,5000,2500,2500,4500,3000,6000,5500,5000,5000,4000,3000,4500,2500,5625,3000,5000,2500,2500,5000,2500,6000,2500,2500,3000,3000,5000,2500,5000,2500,2769,3000,2500,4000,6000,2500)
$colO = @("Región Metropolitana","Región Metropolitana","Región Metropolitana","Región Metropolitana","Región Metropolitana","Provincia de Cautín","Región Metropolitana","Región Metropolitana","Región Metropolitana","Región Metropolitana","Región Metropolitana","Región Metropolitana","Región Metropolitana","Provincia de Cautín","Región Metropolitana","Región Metropolitana","Región Metropolitana","Región Metropolitana","Provincia de Cautín","Provincia de Cautín","Región Metropolitana","Región Metropolitana","Región Metropolitana","Región Metropolitana","Provincia de Cautín","Región Metropolitana","Región Metropolitana","Provincia de Cautín","Región Metropolitana","Región Metropolitana","Región Metropolitana","Provincia de Cautín","Región Metropolitana","Región Metropolitana","Región Metropolitana","Región Metropolitana","Región Metropolitana","Provincia de Cautín","Provincia de Cautín","Provincia de Cautín","Provincia de Cautín","Región Metropolitana","Región Metropolitana","Región Metropolitana","Región Metropolitana","Provincia de Cautín","Región Metropolitana","Región Metropolitana","Región Metropolitana","Región Metropolitana","Región Metropolitana","Región Metropolitana","Provincia de Cautín","Región Metropolitana","Región Metropolitana","Región Metropolitana","Región Metropolitana","Provincia de Cautín","Región Metropolitana","Provincia de Cautín","Región Metropolitana","Región Metropolitana","Región Metropolitana","Región Metropolitana","Región Metropolitana","Provincia de Cautín","Región Metropolitana")
$colP = @(917,1500,833,1000,833,1833,1593,1000,833,1500,833,1167,1500,1917,1000,1000,1000,833,2000,1667,1000,833,1333,1333,1667,1000,1000,2000,833,1000,833,2000,1667,833,833,1500,1000,2000,1833,1667,1667,1333,1000,1500,833,1875,1000,1667,833,833,1667,833,2000,833,833,1000,1000,1667,833,1667,833,923,1000,833,1333,2000,833)

# Style/number-format to apply to the date column of the brand-new row (166),
# copied from the row directly above it so the date renders the same way.
$dateFormat = $ws.Cells.Item(165, 4).NumberFormat

for ($i = 0; $i -lt $colD.Length; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 4).Value2 = $colD[$i]
    $ws.Cells.Item($r, 10).Value2 = $colJ[$i]
    $ws.Cells.Item($r, 11).Value2 = $colK[$i]
    $ws.Cells.Item($r, 12).Value2 = $colL[$i]
    $ws.Cells.Item($r, 13).Value2 = $colM[$i]
    $ws.Cells.Item($r, 15).Value2 = $colO[$i]
    $ws.Cells.Item($r, 16).Value2 = $colP[$i]
}

# Row 166 is a brand new row; fill in the columns that are constant across
# the whole table (they were not previously populated for this row).
$ws.Cells.Item(166, 1).Value2 = $ws.Cells.Item(165, 1).Value2
$ws.Cells.Item(166, 2).Value2 = $ws.Cells.Item(165, 2).Value2
$ws.Cells.Item(166, 3).Value2 = $ws.Cells.Item(165, 3).Value2
$ws.Cells.Item(166, 5).Value2 = $ws.Cells.Item(165, 5).Value2
$ws.Cells.Item(166, 6).Value2 = $ws.Cells.Item(165, 6).Value2
$ws.Cells.Item(166, 7).Value2 = $ws.Cells.Item(165, 7).Value2
$ws.Cells.Item(166, 8).Value2 = $ws.Cells.Item(165, 8).Value2
$ws.Cells.Item(166, 9).Value2 = $ws.Cells.Item(165, 9).Value2
$ws.Cells.Item(166, 14).Value2 = $ws.Cells.Item(165, 14).Value2
$ws.Cells.Item(166, 17).Value2 = $ws.Cells.Item(165, 17).Value2
$ws.Cells.Item(166, 18).Value2 = $ws.Cells.Item(165, 18).Value2

# Match the date-cell style/number format for the newly created row
$ws.Cells.Item(166, 4).NumberFormat = $dateFormat

Write-Host "Dimension now:" $ws.UsedRange.Address()
